$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35; this shifts the existing rows 35..126 down to 36..127
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with the new weekly price record
$ws.Cells.Item(35, 1).Value = 4
$ws.Cells.Item(35, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(35, 3).Value = "Los Lagos"
$ws.Cells.Item(35, 4).Value = Get-Date -Year 2023 -Month 2 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(35, 5).Value = 10
$ws.Cells.Item(35, 6).Value = 100112031
$ws.Cells.Item(35, 7).Value = "Poroto verde"
$ws.Cells.Item(35, 8).Value = "Magnum"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 20
$ws.Cells.Item(35, 11).Value = 27000
$ws.Cells.Item(35, 12).Value = 27000
$ws.Cells.Item(35, 13).Value = 27000
$ws.Cells.Item(35, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(35, 15).Value = "Región Metropolitana"
$ws.Cells.Item(35, 16).Value = 1080
$ws.Cells.Item(35, 17).Value = 25
$ws.Cells.Item(35, 18).Value = "Hortaliza"
